$wb = $excel.ActiveWorkbook

# Rename "Sheet1" worksheet to "Data".
$dataSheet = $wb.Worksheets.Item("Sheet1")
$dataSheet.Name = "Data"
$dataSheet.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 670
$win.ScrollColumn = 1

$dataSheet.Range("J698").Select()
